$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark that currently sits
#     after "12. Affectez les ressources aux taches (Annexe 4)." ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: re-create the "_GoBack" bookmark right after the
#     "...et rectifiez si besoin." paragraph.
#     A zero-length Range sitting exactly at "end of paragraph text,
#     right before the paragraph mark" cannot be used directly with
#     Bookmarks.Add, so insert a temporary marker character, wrap the
#     bookmark around it, then delete the marker again - the bookmark
#     collapses back down to the correct (zero-length) position. ---
$findRng = $d.Content
$found = $findRng.Find.Execute("rectifiez si besoin.", $true, $false, $false,
                                $false, $false, $true, 1, $false, "", 0)
$pos = $findRng.End

$insRng = $d.Range($pos, $pos)
$insRng.InsertAfter("@")
$markRng = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $markRng)
$markRng2 = $d.Range($pos, $pos + 1)
$markRng2.Text = ""

# --- Step 3: update the project start date, 21/03/11 -> 21/03/23 ---
$dateFindRng = $d.Content
$foundDate = $dateFindRng.Find.Execute("21/03/11", $true, $false, $false,
                                        $false, $false, $true, 1, $false,
                                        "", 0)
$dateEnd = $dateFindRng.End
$lastTwo = $d.Range($dateEnd - 2, $dateEnd)
$lastTwo.Text = "23"
